# Clean and sort out translators (sbml-qual part done)
# Update the "Regulator Name" values in column A for specific rows so that
# the interaction list entries line up with their correct regulator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = "MTORC2"
    3  = "PDK1"
    8  = "STAT5"
    9  = "NFAT"
    10 = "AP1"
    12 = "FOXP3"
    17 = "MTOR_DD"
    18 = "NFAT"
    19 = "SMAD3"
    20 = "STAT5"
    21 = "NFAT"
    22 = "NFKAPPAB"
    23 = "AP1"
    24 = "IL2"
    26 = "IL2_EX"
    27 = "IL2"
    28 = "CD25"
    29 = "CD122"
    30 = "CD132"
    52 = "TCR"
    53 = "CD28"
    54 = "PI3K"
    55 = "IL2_EX"
    56 = "IL2R"
    59 = "MTORC2"
    61 = "CD28"
    63 = "PTEN"
    64 = "FOXP3"
    67 = "TCR"
    68 = "CD28"
    69 = "RAS"
    70 = "IL2_EX"
    71 = "IL2R"
}

foreach ($row in $newValues.Keys) {
    $ws.Range("A$row").Value = $newValues[$row]
}
